$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecasts")

$ws.Range("B3").Value = 0.6487449759341348
$ws.Range("E3").Value = 0.8407226460642582

$ws.Range("B4").Value = 0.09976210878375116
$ws.Range("E4").Value = 0.9987624131073437

$ws.Range("B5").Value = 6.472016686205318
$ws.Range("E5").Value = 83.96821787371297

$ws.Range("B6").Value = 0.9870512
$ws.Range("C6").Value = 0.0129488
$ws.Range("E6").Value = 0.4792927
$ws.Range("F6").Value = 0.5207073

$ws.Range("B7").Value = 33.7956664
$ws.Range("C7").Value = 16.1666372
$ws.Range("E7").Value = 24.038751
$ws.Range("F7").Value = 24.8354808
